# This script re-applies a data re-scrape that reordered a handful of
# fixture rows which share the exact same Date/Time within the sheet.
# For each affected group of rows, the entire row's content in columns
# B..AC (id, Div, Div Original Name, Date, HomeTeam, AwayTeam, FTHG, FTAG,
# FTR, odds, ...) is moved to a different row in the group, while column A
# (the static positional row counter) stays put on its own row.
#
# Groups (by worksheet row number):
#   - Rows 5, 6                  : swap
#   - Rows 9, 10, 11, 12         : cyclic rotation (9<-10<-11<-12<-9)
#   - Rows 22, 23                : swap
#   - Rows 41, 42                : swap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC inclusive (id ... PL_AhUnder), i.e. everything except
# column A (the static row counter).
$firstCol = 2   # B
$lastCol  = 29  # AC

function Get-RowValues($row) {
    $rng = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    return $rng.Value()
}

function Set-RowValues($row, $values) {
    $rng = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    $rng.Value = $values
}

# ---- Snapshot current (pre-edit) row contents for every affected row ----
$row5  = Get-RowValues 5
$row6  = Get-RowValues 6

$row9  = Get-RowValues 9
$row10 = Get-RowValues 10
$row11 = Get-RowValues 11
$row12 = Get-RowValues 12

$row22 = Get-RowValues 22
$row23 = Get-RowValues 23

$row41 = Get-RowValues 41
$row42 = Get-RowValues 42

# ---- Apply swap: rows 5 <-> 6 ----
Set-RowValues 5 $row6
Set-RowValues 6 $row5

# ---- Apply cyclic rotation: new9=old10, new10=old11, new11=old12, new12=old9 ----
Set-RowValues 9  $row10
Set-RowValues 10 $row11
Set-RowValues 11 $row12
Set-RowValues 12 $row9

# ---- Apply swap: rows 22 <-> 23 ----
Set-RowValues 22 $row23
Set-RowValues 23 $row22

# ---- Apply swap: rows 41 <-> 42 ----
Set-RowValues 41 $row42
Set-RowValues 42 $row41
